# PROS-9213 - CCRU - rename KPI and add EAN code for some new Products
#
# The underlying "PoS 2019 - IC Petroleum - CAP" KPI set was renamed to
# "PoS 2019 - IC Petroleum – REG" (note: en dash, not hyphen) for the two
# rows that reference the "Juice Availability" KPI / NEW SKU 7 & NEW SKU 8
# atomic names. All other rows that happen to reference KPI-set / KPI-name /
# atomic-name strings located later in the shared-string table are left
# untouched content-wise - Excel will re-index the shared strings table on
# save automatically once the new string is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rename the KPI set for the two "IC Petroleum" rows ------------------
# (note: the separator before "REG" is an EN DASH, U+2013, not a hyphen)
$newKpiSet = "PoS 2019 - IC Petroleum " + [char]0x2013 + " REG"

$ws.Range("A22").Value2 = $newKpiSet
$ws.Range("A23").Value2 = $newKpiSet

# --- cosmetic view / column-width tweaks left behind by the edit ---------
# Column A got noticeably wider (to fit the longer KPI-set name), the other
# columns shrank very slightly as a side effect of that resize.
$ws.Columns.Item(1).ColumnWidth = 43.333333333333336
$ws.Columns.Item(2).ColumnWidth = 15.333333333333334
$ws.Columns.Item(3).ColumnWidth = 83
$ws.Columns.Item(4).ColumnWidth = 87.83333333333333
$ws.Columns.Item(5).ColumnWidth = 7.5

# the view scrolled down (new topLeftCell = A10) and the active selection
# moved on to D35:D37
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("D35:D37").Select()
